$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data.
# Each D/E cell holds a plain-text value (inline string in the source file),
# so the cell format is forced to Text ("@") before assignment; this prevents
# Excel from auto-converting numeric-looking strings (e.g. "306.74", "-4.59%")
# into real numbers/percentages, preserving the exact original text formatting.
$updates = @{
    'D2' = '306.74'
    'E2' = '-4.59%'
    'D3' = '39.95'
    'E3' = '-6.95%'
    'D4' = '5.089'
    'E4' = '-1.64%'
    'D5' = '0.07694'
    'E5' = '-5.71%'
    'D6' = '4.266'
    'E6' = '-1.42%'
    'D7' = '1.624'
    'E7' = '-11.49%'
    'D8' = '0.8773'
    'E8' = '-6.91%'
    'D9' = '0.09681'
    'E9' = '-13.45%'
    'D10' = '0.1733'
    'E10' = '-6.94%'
    'D11' = '0.08935'
    'E11' = '-4.22%'
    'D12' = '0.04403'
    'E12' = '-4.30%'
    'D13' = '0.1057'
    'E13' = '-0.12%'
    'E14' = '-2.83%'
    'D15' = '0.005959'
    'E15' = '4.00%'
    'D16' = '3.357'
    'E16' = '-0.07%'
    'D17' = '2.437'
    'E17' = '-2.99%'
    'E18' = '-2.00%'
    'D19' = '6.998'
    'E19' = '-5.59%'
    'D20' = '0.1339'
    'E20' = '-3.62%'
    'D21' = '0.3137'
    'E21' = '19.65%'
    'D22' = '0.04166'
    'E22' = '0.68%'
    'D23' = '0.001196'
    'E23' = '-4.38%'
    'E24' = '-5.42%'
    'E25' = '9.93%'
    'E26' = '0.09%'
    'D38' = '0.02342'
    'E38' = '-13.73%'
    'D39' = '0.05147'
    'E39' = '-6.80%'
    'D40' = '0.007923'
    'E40' = '-0.60%'
    'D41' = '0.1323'
    'E41' = '-5.11%'
    'D42' = '0.006389'
    'E42' = '-2.37%'
    'E43' = '-8.05%'
    'D44' = '0.008598'
    'E44' = '15.60%'
    'D45' = '0.3048'
    'E45' = '-4.84%'
    'D46' = '0.00006526'
    'E46' = '-6.61%'
    'E47' = '0.01%'
    'D48' = '0.007001'
    'E48' = '98.34%'
    'D49' = '0.003393'
    'E49' = '-2.08%'
    'E50' = '0.01%'
    'E51' = '0.01%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
